# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N ("Late") on the
# "Repayment schedule" sheet, pushing the existing N/O/P columns
# (Late / heading / Outstanding) one column to the right, and leave the
# sheet selection on the newly-shaped table (matches the tabSelected /
# selection state recorded for this sheet after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at N (14th column); existing N->O, O->P, P->Q.
$ws.Columns.Item(14).Insert()

# Match the column width Excel gives the freshly inserted column (copied
# from its left neighbour, column M).
$ws.Columns.Item(14).ColumnWidth = 10.17

# Make "Repayment schedule" the active sheet/tab and move the selection
# to the new bottom-right corner of the table, as recorded in the saved
# view state.
$ws.Activate() | Out-Null
$ws.Range("Q10").Select() | Out-Null
